$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new tracker entry (row 29) for LeetCode 3Sum, cloning the
# formatting of the row directly above it (row 28) so that the new
# cells reuse the existing style definitions instead of creating new ones.
$ws.Rows("28:28").Copy()
$ws.Rows("29:29").Insert(-4121)  # xlShiftDown

$ws.Range("A29").Value = 15
$ws.Range("B29").Value = "3Sum"
$ws.Range("C29").Value = "Medium"
$ws.Range("D29").Value = "Arrays,Two Pointers ,Sorting"
$ws.Range("E29").Value = 45699

# Update the view state to match where the workbook was scrolled/selected
# when the row was added.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
